$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 50520.99
$ws.Range("D2").Value = 141960.2
$ws.Range("E2").Value = 164490.68
$ws.Range("F2").Value = 327246.8
$ws.Range("G2").Value = 402851.09
$ws.Range("C3").Value = 20344.44
$ws.Range("D3").Value = 150440.59
$ws.Range("E3").Value = 189967.91
$ws.Range("F3").Value = 347725.48
$ws.Range("G3").Value = 457315.71
$ws.Range("D4").Value = 167917.14
$ws.Range("E4").Value = 196889.99
$ws.Range("F4").Value = 325429.09
$ws.Range("G4").Value = 362961.66
$ws.Range("C5").Value = 44518.43
$ws.Range("D5").Value = 130937.67
$ws.Range("E5").Value = 220983.18
$ws.Range("F5").Value = 411174.53
$ws.Range("G5").Value = 415905.04
$ws.Range("D6").Value = 116057.63
$ws.Range("E6").Value = 222492.71
$ws.Range("F6").Value = 386420.95
$ws.Range("G6").Value = 134984.25
$ws.Range("C7").Value = 64428.92
$ws.Range("D7").Value = 142478.65
$ws.Range("E7").Value = 249067.64
$ws.Range("F7").Value = 412060.64
$ws.Range("C8").Value = 72571.03
$ws.Range("D8").Value = 139276.88
$ws.Range("E8").Value = 224412.72
$ws.Range("F8").Value = 395421.28
$ws.Range("C9").Value = 116521.88
$ws.Range("D9").Value = 152291.48
$ws.Range("E9").Value = 272902.32
$ws.Range("F9").Value = 390627.44
$ws.Range("B10").Value = 40596.4
$ws.Range("D10").Value = 157483.4
$ws.Range("E10").Value = 301536.61
$ws.Range("F10").Value = 396072.84
$ws.Range("C11").Value = 119518.98
$ws.Range("D11").Value = 134933.1
$ws.Range("E11").Value = 335641.61
$ws.Range("F11").Value = 434296.15
$ws.Range("B12").Value = 18815.61
$ws.Range("C12").Value = 129663.99
$ws.Range("D12").Value = 151498.73
$ws.Range("E12").Value = 253267.57
$ws.Range("F12").Value = 309659.47
$ws.Range("B13").Value = 23098.81
$ws.Range("C13").Value = 107309.77
$ws.Range("E13").Value = 254321.08
$ws.Range("F13").Value = 381298.1
